# Corrects the zip-code range for Rondonia (row 26, dLocal code 22) on Hoja1,
# and records the page-setup / selection state present in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix the min/max zip-code range values for Rondonia
$ws.Range("D26").Value = 76800000
$ws.Range("E26").Value = 76999999

# Page setup: Letter/A4-class paper, portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the view/selection state saved with the workbook
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F24").Select()
